$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '24.729.56'
$ws.Range("E2").Value = '  +0.75%  '

$ws.Range("D3").Value = '1.700.64'
$ws.Range("E3").Value = '  +0.45%  '

$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").Value = '315.93'
$ws.Range("E5").Value = '  +0.13%  '

$ws.Range("D6").Value = '1.002'
$ws.Range("E6").Value = '  +0.00%  '

$ws.Range("D7").Value = '0.3929'
$ws.Range("E7").Value = '  -0.07%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4050'
$ws.Range("E8").Value = '  +1.22%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.520'
$ws.Range("E9").Value = '  -0.08%  '

$ws.Range("E10").Value = '  +0.07%  '

$ws.Range("D11").Value = '52.65'
$ws.Range("E11").Value = '  -1.10%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08860'
$ws.Range("E12").Value = '  +1.56%  '

$ws.Range("D13").Value = '7.451'
$ws.Range("E13").Value = '  +3.73%  '

$ws.Range("E14").Value = '  +2.52%  '

$ws.Range("D15").Value = '8.114'
$ws.Range("E15").Value = '  +6.99%  '

$ws.Range("E16").Value = '  +0.51%  '

$ws.Range("D17").Value = '1.700.48'
$ws.Range("E17").Value = '  -0.03%  '

$ws.Range("D18").Value = '99.44'
$ws.Range("E18").Value = '  -0.23%  '

$ws.Range("D19").Value = '0.07058'
$ws.Range("E19").Value = '  +0.19%  '

$ws.Range("E20").Value = '  +1.30%  '

$ws.Range("D21").Value = '7.076'
$ws.Range("E21").Value = '  +3.56%  '

$ws.Range("D22").Value = '1.006'
$ws.Range("E22").Value = '  +0.44%  '

$ws.Range("D23").Value = '14.74'
$ws.Range("E23").Value = '  +5.13%  '

$ws.Range("D24").Value = '24.731.83'

$ws.Range("D25").Value = '3.152'
$ws.Range("E25").Value = '  +5.28%  '

$ws.Range("E26").Value = '  +1.40%  '

$ws.Range("D27").Value = '22.67'
$ws.Range("E27").Value = '  +1.65%  '

$ws.Range("D28").Value = '164.53'
$ws.Range("E28").Value = '  +2.39%  '

$ws.Range("D29").Value = '9.011'
$ws.Range("E29").Value = '  +20.67%  '

$ws.Range("D30").Value = '135.71'
$ws.Range("E30").Value = '  +1.11%  '

$ws.Range("D31").Value = '5.148'
$ws.Range("E31").Value = '  -1.20%  '

$ws.Range("D32").Value = '7.737'
$ws.Range("E32").Value = '  +6.43%  '

$ws.Range("D33").Value = '0.09051'
$ws.Range("E33").Value = '  +6.31%  '

$ws.Range("D34").Value = '1.072'
$ws.Range("E34").Value = '  -1.71%  '

$ws.Range("D35").Value = '0.02948'
$ws.Range("E35").Value = '  +7.36%  '

# Row 36/37 swap: WEMIXTOKEN <-> Algorand
$ws.Range("B36").Value = 'Algorand'
$ws.Range("C36").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D36").Value = '0.2765'
$ws.Range("E36").Value = '  +2.15%  '

$ws.Range("B37").Value = 'WEMIXTOKEN'
$ws.Range("C37").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D37").Value = '1.962'
$ws.Range("E37").Value = '  +0.36%  '

# Row 44/45 swap: NEARProtocol <-> Decentraland
$ws.Range("B44").Value = 'Decentraland'
$ws.Range("C44").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D44").Value = '0.7203'
$ws.Range("E44").Value = '  +0.45%  '

$ws.Range("B45").Value = 'NEARProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D45").Value = '2.597'
$ws.Range("E45").Value = '  +3.03%  '

$ws.Range("D38").Value = '11.05'
$ws.Range("E38").Value = '  -2.64%  '

$ws.Range("E39").Value = '  +0.52%  '

$ws.Range("D40").Value = '0.09242'
$ws.Range("E40").Value = '  +2.46%  '

$ws.Range("D41").Value = '1.468'
$ws.Range("E41").Value = '  -0.51%  '

$ws.Range("D42").Value = '0.7747'
$ws.Range("E42").Value = '  +1.66%  '

$ws.Range("D43").Value = '16.25'
$ws.Range("E43").Value = '  +6.11%  '

$ws.Range("D46").Value = '4.206'
$ws.Range("E46").Value = '  +0.06%  '

$ws.Range("E47").Value = '  +4.05%  '

$ws.Range("E48").Value = '  +0.08%  '

$ws.Range("D49").Value = '139.91'
$ws.Range("E49").Value = '  -0.65%  '

$ws.Range("D50").Value = '0.07975'
$ws.Range("E50").Value = '  -0.27%  '

$ws.Range("D51").Value = '89.71'
$ws.Range("E51").Value = '  +2.24%  '
